$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value2 = 'FAPs'
$ws.Cells.Item(2,2).Value2 = 'Qrfp'
$ws.Cells.Item(2,3).Value2 = 'P2ry14'
$ws.Cells.Item(2,4).Value2 = 'ECs'
$ws.Cells.Item(2,5).Value2 = 2
$ws.Cells.Item(2,6).Value2 = 0.6666666666666666
$ws.Cells.Item(2,7).Value2 = 0.5328176666666667
$ws.Cells.Item(2,8).Value2 = 1.598453
$ws.Cells.Item(2,9).Value2 = 0.2001388429706536
$ws.Cells.Item(2,10).Value2 = 0.213328135545356
$ws.Cells.Item(2,11).Value2 = 2
$ws.Cells.Item(2,12).Value2 = 0.6666666666666666
$ws.Cells.Item(2,13).Value2 = 0.107352
$ws.Cells.Item(2,14).Value2 = 0.322056
$ws.Cells.Item(2,15).Value2 = 0.004186491276943991
$ws.Cells.Item(2,16).Value2 = 0.004195990443300053
$ws.Cells.Item(2,17).Value2 = 0.05719904215200001
$ws.Cells.Item(2,18).Value2 = 0.514791379368
$ws.Cells.Item(2,19).Value2 = 0.0008378795202743046
$ws.Cells.Item(2,20).Value2 = 0.0008951228180353322

$ws.Cells.Item(3,1).Value2 = 'FAPs'
$ws.Cells.Item(3,2).Value2 = 'Qrfp'
$ws.Cells.Item(3,3).Value2 = 'P2ry14'
$ws.Cells.Item(3,4).Value2 = 'FAPs'
$ws.Cells.Item(3,5).Value2 = 2
$ws.Cells.Item(3,6).Value2 = 0.6666666666666666
$ws.Cells.Item(3,7).Value2 = 0.5328176666666667
$ws.Cells.Item(3,8).Value2 = 1.598453
$ws.Cells.Item(3,9).Value2 = 0.2001388429706536
$ws.Cells.Item(3,10).Value2 = 0.213328135545356
$ws.Cells.Item(3,11).Value2 = 3
$ws.Cells.Item(3,12).Value2 = 1
$ws.Cells.Item(3,13).Value2 = 1.843770666666667
$ws.Cells.Item(3,14).Value2 = 5.531312
$ws.Cells.Item(3,15).Value2 = 0.07190299028136603
$ws.Cells.Item(3,16).Value2 = 0.07206613846943047
$ws.Cells.Item(3,17).Value2 = 0.9823935844817779
$ws.Cells.Item(3,18).Value2 = 8.841542260336
$ws.Cells.Item(3,19).Value2 = 0.01439058128104275
$ws.Cells.Item(3,20).Value2 = 0.01537373495563706

$ws.Cells.Item(4,1).Value2 = 'FAPs'
$ws.Cells.Item(4,2).Value2 = 'Qrfp'
$ws.Cells.Item(4,3).Value2 = 'P2ry14'
$ws.Cells.Item(4,4).Value2 = 'Inflammatory-Mac'
$ws.Cells.Item(4,5).Value2 = 2
$ws.Cells.Item(4,6).Value2 = 0.6666666666666666
$ws.Cells.Item(4,7).Value2 = 0.5328176666666667
$ws.Cells.Item(4,8).Value2 = 1.598453
$ws.Cells.Item(4,9).Value2 = 0.2001388429706536
$ws.Cells.Item(4,10).Value2 = 0.213328135545356
$ws.Cells.Item(4,11).Value2 = 3
$ws.Cells.Item(4,12).Value2 = 1
$ws.Cells.Item(4,13).Value2 = 13.883856
$ws.Cells.Item(4,14).Value2 = 41.651568
$ws.Cells.Item(4,15).Value2 = 0.5414397685590067
$ws.Cells.Item(4,16).Value2 = 0.5426682976763739
$ws.Cells.Item(4,17).Value2 = 7.397563758256001
$ws.Cells.Item(4,18).Value2 = 66.578073824304
$ws.Cells.Item(4,19).Value2 = 0.1083631288176981
$ws.Cells.Item(4,20).Value2 = 0.1157664161628731

$ws.Cells.Item(5,1).Value2 = 'FAPs'
$ws.Cells.Item(5,2).Value2 = 'Qrfp'
$ws.Cells.Item(5,3).Value2 = 'P2ry14'
$ws.Cells.Item(5,4).Value2 = 'MuSCs'
$ws.Cells.Item(5,5).Value2 = 2
$ws.Cells.Item(5,6).Value2 = 0.6666666666666666
$ws.Cells.Item(5,7).Value2 = 0.5328176666666667
$ws.Cells.Item(5,8).Value2 = 1.598453
$ws.Cells.Item(5,9).Value2 = 0.2001388429706536
$ws.Cells.Item(5,10).Value2 = 0.213328135545356
$ws.Cells.Item(5,11).Value2 = 2
$ws.Cells.Item(5,12).Value2 = 1
$ws.Cells.Item(5,13).Value2 = 0.1741535
$ws.Cells.Item(5,14).Value2 = 0.348307
$ws.Cells.Item(5,15).Value2 = 0.006791602472234009
$ws.Cells.Item(5,16).Value2 = 0.00453800843124957
$ws.Cells.Item(5,17).Value2 = 0.09279206151183336
$ws.Cells.Item(5,18).Value2 = 0.5567523690710001
$ws.Cells.Item(5,19).Value2 = 0.001359263460709545
$ws.Cells.Item(5,20).Value2 = 0.0009680848777275767

$ws.Cells.Item(6,1).Value2 = 'FAPs'
$ws.Cells.Item(6,2).Value2 = 'Qrfp'
$ws.Cells.Item(6,3).Value2 = 'P2ry14'
$ws.Cells.Item(6,4).Value2 = 'Resolving-Mac'
$ws.Cells.Item(6,5).Value2 = 2
$ws.Cells.Item(6,6).Value2 = 0.6666666666666666
$ws.Cells.Item(6,7).Value2 = 0.5328176666666667
$ws.Cells.Item(6,8).Value2 = 1.598453
$ws.Cells.Item(6,9).Value2 = 0.2001388429706536
$ws.Cells.Item(6,10).Value2 = 0.213328135545356
$ws.Cells.Item(6,11).Value2 = 3
$ws.Cells.Item(6,12).Value2 = 1
$ws.Cells.Item(6,13).Value2 = 9.633343333333334
$ws.Cells.Item(6,14).Value2 = 28.90003
$ws.Cells.Item(6,15).Value2 = 0.3756791474104492
$ws.Cells.Item(6,16).Value2 = 0.3765315649796458
$ws.Cells.Item(6,17).Value2 = 5.132815517065556
$ws.Cells.Item(6,18).Value2 = 46.19533965359
$ws.Cells.Item(6,19).Value2 = 0.07518798989092894
$ws.Cells.Item(6,20).Value2 = 0.08032477673108292

$ws.Cells.Item(7,1).Value2 = 'Inflammatory-Mac'
$ws.Cells.Item(7,2).Value2 = 'Qrfp'
$ws.Cells.Item(7,3).Value2 = 'P2ry14'
$ws.Cells.Item(7,4).Value2 = 'ECs'
$ws.Cells.Item(7,5).Value2 = 3
$ws.Cells.Item(7,6).Value2 = 1
$ws.Cells.Item(7,7).Value2 = 0.8624390000000001
$ws.Cells.Item(7,8).Value2 = 2.587317
$ws.Cells.Item(7,9).Value2 = 0.323952365680006
$ws.Cells.Item(7,10).Value2 = 0.3453010577569712
$ws.Cells.Item(7,11).Value2 = 2
$ws.Cells.Item(7,12).Value2 = 0.6666666666666666
$ws.Cells.Item(7,13).Value2 = 0.107352
$ws.Cells.Item(7,14).Value2 = 0.322056
$ws.Cells.Item(7,15).Value2 = 0.004186491276943991
$ws.Cells.Item(7,16).Value2 = 0.004195990443300053
$ws.Cells.Item(7,17).Value2 = 0.09258455152800001
$ws.Cells.Item(7,18).Value2 = 0.8332609637520001
$ws.Cells.Item(7,19).Value2 = 0.001356223753064715
$ws.Cells.Item(7,20).Value2 = 0.001448879938409651

$ws.Cells.Item(8,1).Value2 = 'Inflammatory-Mac'
$ws.Cells.Item(8,2).Value2 = 'Qrfp'
$ws.Cells.Item(8,3).Value2 = 'P2ry14'
$ws.Cells.Item(8,4).Value2 = 'FAPs'
$ws.Cells.Item(8,5).Value2 = 3
$ws.Cells.Item(8,6).Value2 = 1
$ws.Cells.Item(8,7).Value2 = 0.8624390000000001
$ws.Cells.Item(8,8).Value2 = 2.587317
$ws.Cells.Item(8,9).Value2 = 0.323952365680006
$ws.Cells.Item(8,10).Value2 = 0.3453010577569712
$ws.Cells.Item(8,11).Value2 = 3
$ws.Cells.Item(8,12).Value2 = 1
$ws.Cells.Item(8,13).Value2 = 1.843770666666667
$ws.Cells.Item(8,14).Value2 = 5.531312
$ws.Cells.Item(8,15).Value2 = 0.07190299028136603
$ws.Cells.Item(8,16).Value2 = 0.07206613846943047
$ws.Cells.Item(8,17).Value2 = 1.590139729989333
$ws.Cells.Item(8,18).Value2 = 14.311257569904
$ws.Cells.Item(8,19).Value2 = 0.023293143801115
$ws.Cells.Item(8,20).Value2 = 0.02488451384195469

$ws.Cells.Item(9,1).Value2 = 'Inflammatory-Mac'
$ws.Cells.Item(9,2).Value2 = 'Qrfp'
$ws.Cells.Item(9,3).Value2 = 'P2ry14'
$ws.Cells.Item(9,4).Value2 = 'Inflammatory-Mac'
$ws.Cells.Item(9,5).Value2 = 3
$ws.Cells.Item(9,6).Value2 = 1
$ws.Cells.Item(9,7).Value2 = 0.8624390000000001
$ws.Cells.Item(9,8).Value2 = 2.587317
$ws.Cells.Item(9,9).Value2 = 0.323952365680006
$ws.Cells.Item(9,10).Value2 = 0.3453010577569712
$ws.Cells.Item(9,11).Value2 = 3
$ws.Cells.Item(9,12).Value2 = 1
$ws.Cells.Item(9,13).Value2 = 13.883856
$ws.Cells.Item(9,14).Value2 = 41.651568
$ws.Cells.Item(9,15).Value2 = 0.5414397685590067
$ws.Cells.Item(9,16).Value2 = 0.5426682976763739
$ws.Cells.Item(9,17).Value2 = 11.973978884784
$ws.Cells.Item(9,18).Value2 = 107.765809963056
$ws.Cells.Item(9,19).Value2 = 0.1754006938979251
$ws.Cells.Item(9,20).Value2 = 0.1873839371988268

$ws.Cells.Item(10,1).Value2 = 'Inflammatory-Mac'
$ws.Cells.Item(10,2).Value2 = 'Qrfp'
$ws.Cells.Item(10,3).Value2 = 'P2ry14'
$ws.Cells.Item(10,4).Value2 = 'MuSCs'
$ws.Cells.Item(10,5).Value2 = 3
$ws.Cells.Item(10,6).Value2 = 1
$ws.Cells.Item(10,7).Value2 = 0.8624390000000001
$ws.Cells.Item(10,8).Value2 = 2.587317
$ws.Cells.Item(10,9).Value2 = 0.323952365680006
$ws.Cells.Item(10,10).Value2 = 0.3453010577569712
$ws.Cells.Item(10,11).Value2 = 2
$ws.Cells.Item(10,12).Value2 = 1
$ws.Cells.Item(10,13).Value2 = 0.1741535
$ws.Cells.Item(10,14).Value2 = 0.348307
$ws.Cells.Item(10,15).Value2 = 0.006791602472234009
$ws.Cells.Item(10,16).Value2 = 0.00453800843124957
$ws.Cells.Item(10,17).Value2 = 0.1501967703865
$ws.Cells.Item(10,18).Value2 = 0.9011806223190001
$ws.Cells.Item(10,19).Value2 = 0.002200155687638384
$ws.Cells.Item(10,20).Value2 = 0.00156697911142053

$ws.Cells.Item(11,1).Value2 = 'Inflammatory-Mac'
$ws.Cells.Item(11,2).Value2 = 'Qrfp'
$ws.Cells.Item(11,3).Value2 = 'P2ry14'
$ws.Cells.Item(11,4).Value2 = 'Resolving-Mac'
$ws.Cells.Item(11,5).Value2 = 3
$ws.Cells.Item(11,6).Value2 = 1
$ws.Cells.Item(11,7).Value2 = 0.8624390000000001
$ws.Cells.Item(11,8).Value2 = 2.587317
$ws.Cells.Item(11,9).Value2 = 0.323952365680006
$ws.Cells.Item(11,10).Value2 = 0.3453010577569712
$ws.Cells.Item(11,11).Value2 = 3
$ws.Cells.Item(11,12).Value2 = 1
$ws.Cells.Item(11,13).Value2 = 9.633343333333334
$ws.Cells.Item(11,14).Value2 = 28.90003
$ws.Cells.Item(11,15).Value2 = 0.3756791474104492
$ws.Cells.Item(11,16).Value2 = 0.3765315649796458
$ws.Cells.Item(11,17).Value2 = 8.308170991056668
$ws.Cells.Item(11,18).Value2 = 74.77353891951
$ws.Cells.Item(11,19).Value2 = 0.1217021485402627
$ws.Cells.Item(11,20).Value2 = 0.1300167476663594

$ws.Cells.Item(12,1).Value2 = 'MuSCs'
$ws.Cells.Item(12,2).Value2 = 'Qrfp'
$ws.Cells.Item(12,3).Value2 = 'P2ry14'
$ws.Cells.Item(12,4).Value2 = 'ECs'
$ws.Cells.Item(12,5).Value2 = 2
$ws.Cells.Item(12,6).Value2 = 1
$ws.Cells.Item(12,7).Value2 = 0.4937895
$ws.Cells.Item(12,8).Value2 = 0.987579
$ws.Cells.Item(12,9).Value2 = 0.1854789459578559
$ws.Cells.Item(12,10).Value2 = 0.1318014272385532
$ws.Cells.Item(12,11).Value2 = 2
$ws.Cells.Item(12,12).Value2 = 0.6666666666666666
$ws.Cells.Item(12,13).Value2 = 0.107352
$ws.Cells.Item(12,14).Value2 = 0.322056
$ws.Cells.Item(12,15).Value2 = 0.004186491276943991
$ws.Cells.Item(12,16).Value2 = 0.004195990443300053
$ws.Cells.Item(12,17).Value2 = 0.053009290404
$ws.Cells.Item(12,18).Value2 = 0.318055742424
$ws.Cells.Item(12,19).Value2 = 0.0007765059893093297
$ws.Cells.Item(12,20).Value2 = 0.0005530375291062766

$ws.Cells.Item(13,1).Value2 = 'MuSCs'
$ws.Cells.Item(13,2).Value2 = 'Qrfp'
$ws.Cells.Item(13,3).Value2 = 'P2ry14'
$ws.Cells.Item(13,4).Value2 = 'FAPs'
$ws.Cells.Item(13,5).Value2 = 2
$ws.Cells.Item(13,6).Value2 = 1
$ws.Cells.Item(13,7).Value2 = 0.4937895
$ws.Cells.Item(13,8).Value2 = 0.987579
$ws.Cells.Item(13,9).Value2 = 0.1854789459578559
$ws.Cells.Item(13,10).Value2 = 0.1318014272385532
$ws.Cells.Item(13,11).Value2 = 3
$ws.Cells.Item(13,12).Value2 = 1
$ws.Cells.Item(13,13).Value2 = 1.843770666666667
$ws.Cells.Item(13,14).Value2 = 5.531312
$ws.Cells.Item(13,15).Value2 = 0.07190299028136603
$ws.Cells.Item(13,16).Value2 = 0.07206613846943047
$ws.Cells.Item(13,17).Value2 = 0.910434595608
$ws.Cells.Item(13,18).Value2 = 5.462607573648
$ws.Cells.Item(13,19).Value2 = 0.01333649084860573
$ws.Cells.Item(13,20).Value2 = 0.009498419905842142

$ws.Cells.Item(14,1).Value2 = 'MuSCs'
$ws.Cells.Item(14,2).Value2 = 'Qrfp'
$ws.Cells.Item(14,3).Value2 = 'P2ry14'
$ws.Cells.Item(14,4).Value2 = 'Inflammatory-Mac'
$ws.Cells.Item(14,5).Value2 = 2
$ws.Cells.Item(14,6).Value2 = 1
$ws.Cells.Item(14,7).Value2 = 0.4937895
$ws.Cells.Item(14,8).Value2 = 0.987579
$ws.Cells.Item(14,9).Value2 = 0.1854789459578559
$ws.Cells.Item(14,10).Value2 = 0.1318014272385532
$ws.Cells.Item(14,11).Value2 = 3
$ws.Cells.Item(14,12).Value2 = 1
$ws.Cells.Item(14,13).Value2 = 13.883856
$ws.Cells.Item(14,14).Value2 = 41.651568
$ws.Cells.Item(14,15).Value2 = 0.5414397685590067
$ws.Cells.Item(14,16).Value2 = 0.5426682976763739
$ws.Cells.Item(14,17).Value2 = 6.855702312311999
$ws.Cells.Item(14,18).Value2 = 41.134213873872
$ws.Cells.Item(14,19).Value2 = 0.10042567757199
$ws.Cells.Item(14,20).Value2 = 0.07152445615086214

$ws.Cells.Item(15,1).Value2 = 'MuSCs'
$ws.Cells.Item(15,2).Value2 = 'Qrfp'
$ws.Cells.Item(15,3).Value2 = 'P2ry14'
$ws.Cells.Item(15,4).Value2 = 'MuSCs'
$ws.Cells.Item(15,5).Value2 = 2
$ws.Cells.Item(15,6).Value2 = 1
$ws.Cells.Item(15,7).Value2 = 0.4937895
$ws.Cells.Item(15,8).Value2 = 0.987579
$ws.Cells.Item(15,9).Value2 = 0.1854789459578559
$ws.Cells.Item(15,10).Value2 = 0.1318014272385532
$ws.Cells.Item(15,11).Value2 = 2
$ws.Cells.Item(15,12).Value2 = 1
$ws.Cells.Item(15,13).Value2 = 0.1741535
$ws.Cells.Item(15,14).Value2 = 0.348307
$ws.Cells.Item(15,15).Value2 = 0.006791602472234009
$ws.Cells.Item(15,16).Value2 = 0.00453800843124957
$ws.Cells.Item(15,17).Value2 = 0.08599516968825001
$ws.Cells.Item(15,18).Value2 = 0.343980678753
$ws.Cells.Item(15,19).Value2 = 0.001259699267914732
$ws.Cells.Item(15,20).Value2 = 0.0005981159880592814

$ws.Cells.Item(16,1).Value2 = 'MuSCs'
$ws.Cells.Item(16,2).Value2 = 'Qrfp'
$ws.Cells.Item(16,3).Value2 = 'P2ry14'
$ws.Cells.Item(16,4).Value2 = 'Resolving-Mac'
$ws.Cells.Item(16,5).Value2 = 2
$ws.Cells.Item(16,6).Value2 = 1
$ws.Cells.Item(16,7).Value2 = 0.4937895
$ws.Cells.Item(16,8).Value2 = 0.987579
$ws.Cells.Item(16,9).Value2 = 0.1854789459578559
$ws.Cells.Item(16,10).Value2 = 0.1318014272385532
$ws.Cells.Item(16,11).Value2 = 3
$ws.Cells.Item(16,12).Value2 = 1
$ws.Cells.Item(16,13).Value2 = 9.633343333333334
$ws.Cells.Item(16,14).Value2 = 28.90003
$ws.Cells.Item(16,15).Value2 = 0.3756791474104492
$ws.Cells.Item(16,16).Value2 = 0.3765315649796458
$ws.Cells.Item(16,17).Value2 = 4.756843787895001
$ws.Cells.Item(16,18).Value2 = 28.54106272737
$ws.Cells.Item(16,19).Value2 = 0.0696805722800361
$ws.Cells.Item(16,20).Value2 = 0.04962739766468337

$ws.Cells.Item(17,1).Value2 = 'Resolving-Mac'
$ws.Cells.Item(17,2).Value2 = 'Qrfp'
$ws.Cells.Item(17,3).Value2 = 'P2ry14'
$ws.Cells.Item(17,4).Value2 = 'ECs'
$ws.Cells.Item(17,5).Value2 = 2
$ws.Cells.Item(17,6).Value2 = 0.6666666666666666
$ws.Cells.Item(17,7).Value2 = 0.773194
$ws.Cells.Item(17,8).Value2 = 2.319582
$ws.Cells.Item(17,9).Value2 = 0.2904298453914845
$ws.Cells.Item(17,10).Value2 = 0.3095693794591195
$ws.Cells.Item(17,11).Value2 = 2
$ws.Cells.Item(17,12).Value2 = 0.6666666666666666
$ws.Cells.Item(17,13).Value2 = 0.107352
$ws.Cells.Item(17,14).Value2 = 0.322056
$ws.Cells.Item(17,15).Value2 = 0.004186491276943991
$ws.Cells.Item(17,16).Value2 = 0.004195990443300053
$ws.Cells.Item(17,17).Value2 = 0.08300392228800001
$ws.Cells.Item(17,18).Value2 = 0.7470353005920001
$ws.Cells.Item(17,19).Value2 = 0.001215882014295642
$ws.Cells.Item(17,20).Value2 = 0.001298950157748793

$ws.Cells.Item(18,1).Value2 = 'Resolving-Mac'
$ws.Cells.Item(18,2).Value2 = 'Qrfp'
$ws.Cells.Item(18,3).Value2 = 'P2ry14'
$ws.Cells.Item(18,4).Value2 = 'FAPs'
$ws.Cells.Item(18,5).Value2 = 2
$ws.Cells.Item(18,6).Value2 = 0.6666666666666666
$ws.Cells.Item(18,7).Value2 = 0.773194
$ws.Cells.Item(18,8).Value2 = 2.319582
$ws.Cells.Item(18,9).Value2 = 0.2904298453914845
$ws.Cells.Item(18,10).Value2 = 0.3095693794591195
$ws.Cells.Item(18,11).Value2 = 3
$ws.Cells.Item(18,12).Value2 = 1
$ws.Cells.Item(18,13).Value2 = 1.843770666666667
$ws.Cells.Item(18,14).Value2 = 5.531312
$ws.Cells.Item(18,15).Value2 = 0.07190299028136603
$ws.Cells.Item(18,16).Value2 = 0.07206613846943047
$ws.Cells.Item(18,17).Value2 = 1.425592416842667
$ws.Cells.Item(18,18).Value2 = 12.830331751584
$ws.Cells.Item(18,19).Value2 = 0.02088277435060255
$ws.Cells.Item(18,20).Value2 = 0.02230946976599657

$ws.Cells.Item(19,1).Value2 = 'Resolving-Mac'
$ws.Cells.Item(19,2).Value2 = 'Qrfp'
$ws.Cells.Item(19,3).Value2 = 'P2ry14'
$ws.Cells.Item(19,4).Value2 = 'Inflammatory-Mac'
$ws.Cells.Item(19,5).Value2 = 2
$ws.Cells.Item(19,6).Value2 = 0.6666666666666666
$ws.Cells.Item(19,7).Value2 = 0.773194
$ws.Cells.Item(19,8).Value2 = 2.319582
$ws.Cells.Item(19,9).Value2 = 0.2904298453914845
$ws.Cells.Item(19,10).Value2 = 0.3095693794591195
$ws.Cells.Item(19,11).Value2 = 3
$ws.Cells.Item(19,12).Value2 = 1
$ws.Cells.Item(19,13).Value2 = 13.883856
$ws.Cells.Item(19,14).Value2 = 41.651568
$ws.Cells.Item(19,15).Value2 = 0.5414397685590067
$ws.Cells.Item(19,16).Value2 = 0.5426682976763739
$ws.Cells.Item(19,17).Value2 = 10.734914156064
$ws.Cells.Item(19,18).Value2 = 96.614227404576
$ws.Cells.Item(19,19).Value2 = 0.1572502682713935
$ws.Cells.Item(19,20).Value2 = 0.1679934881638118

$ws.Cells.Item(20,1).Value2 = 'Resolving-Mac'
$ws.Cells.Item(20,2).Value2 = 'Qrfp'
$ws.Cells.Item(20,3).Value2 = 'P2ry14'
$ws.Cells.Item(20,4).Value2 = 'MuSCs'
$ws.Cells.Item(20,5).Value2 = 2
$ws.Cells.Item(20,6).Value2 = 0.6666666666666666
$ws.Cells.Item(20,7).Value2 = 0.773194
$ws.Cells.Item(20,8).Value2 = 2.319582
$ws.Cells.Item(20,9).Value2 = 0.2904298453914845
$ws.Cells.Item(20,10).Value2 = 0.3095693794591195
$ws.Cells.Item(20,11).Value2 = 2
$ws.Cells.Item(20,12).Value2 = 1
$ws.Cells.Item(20,13).Value2 = 0.1741535
$ws.Cells.Item(20,14).Value2 = 0.348307
$ws.Cells.Item(20,15).Value2 = 0.006791602472234009
$ws.Cells.Item(20,16).Value2 = 0.00453800843124957
$ws.Cells.Item(20,17).Value2 = 0.134654441279
$ws.Cells.Item(20,18).Value2 = 0.8079266476740001
$ws.Cells.Item(20,19).Value2 = 0.001972484055971347
$ws.Cells.Item(20,20).Value2 = 0.001404828454042182

$ws.Cells.Item(21,1).Value2 = 'Resolving-Mac'
$ws.Cells.Item(21,2).Value2 = 'Qrfp'
$ws.Cells.Item(21,3).Value2 = 'P2ry14'
$ws.Cells.Item(21,4).Value2 = 'Resolving-Mac'
$ws.Cells.Item(21,5).Value2 = 2
$ws.Cells.Item(21,6).Value2 = 0.6666666666666666
$ws.Cells.Item(21,7).Value2 = 0.773194
$ws.Cells.Item(21,8).Value2 = 2.319582
$ws.Cells.Item(21,9).Value2 = 0.2904298453914845
$ws.Cells.Item(21,10).Value2 = 0.3095693794591195
$ws.Cells.Item(21,11).Value2 = 3
$ws.Cells.Item(21,12).Value2 = 1
$ws.Cells.Item(21,13).Value2 = 9.633343333333334
$ws.Cells.Item(21,14).Value2 = 28.90003
$ws.Cells.Item(21,15).Value2 = 0.3756791474104492
$ws.Cells.Item(21,16).Value2 = 0.3765315649796458
$ws.Cells.Item(21,17).Value2 = 7.448443265273334
$ws.Cells.Item(21,18).Value2 = 67.03598938746001
$ws.Cells.Item(21,19).Value2 = 0.1091084366992215
$ws.Cells.Item(21,20).Value2 = 0.1165626429175201
